$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header column H (new "comment" column) ---
$ws.Range("H1").Value = "//Debajo de esto se puede "
$ws.Range("H1").HorizontalAlignment = -4108   # xlCenter, matches style of A1:G1 (s="1")

$ws.Range("H2").Value = "//Poner todo lo que quiera"
$ws.Range("H3").Value = "//Solo escriba la // siempre"

# --- Row 5 instructions (entered in this order so the shared-string table
#     indices line up with how the workbook was originally authored) ---
$ws.Range("A5").Value = "//Aquí puede escribir"
$ws.Range("B5").Value = "//Tambien"
$ws.Range("C5").Value = "//Si lo hace"

$ws.Range("E5").Value = "//toda esta hilera"
$ws.Range("F5").Value = "//Y no ponga valores"
$ws.Range("G5").Value = "//Saltados, siempre ponga todo"

$ws.Range("D5").Value = "//rellene"

$ws.Range("H5").Value = "//o dejelo vacio"

# --- Column width tweaks ---
$ws.Columns.Item(6).ColumnWidth = 17.83   # column F a bit wider
$ws.Columns.Item(8).ColumnWidth = 28.33   # new column H width

# --- Selection moves to the newly filled H5 cell ---
$ws.Range("H5").Select() | Out-Null
